$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: lunch moved away from here -> becomes a regular "-" slot
$ws.Range("A8").Value = "11:30"
$ws.Range("B8:F8").Value = "-"

# Row 9: lunch now starts at 12:20
$ws.Range("A9").Value = "12:20"
$ws.Range("B9:F9").Value = "Almoço"

# Row 10: time shifts later by one slot, stays a "-" slot
$ws.Range("A10").Value = "13:00"
$ws.Range("B10:F10").Value = "-"

# Row 11: time shifts later by one slot, stays a "-" slot
$ws.Range("A11").Value = "13:50"
$ws.Range("B11:F11").Value = "-"

# Row 12: break moved away from here -> becomes a regular "-" slot
$ws.Range("A12").Value = "14:40"
$ws.Range("B12:F12").Value = "-"

# Row 13: break now starts at 15:30
$ws.Range("A13").Value = "15:30"
$ws.Range("B13:F13").Value = "Intervalo"

# Row 14 (new): 15:50 slot
$ws.Range("A14").Value = "15:50"
$ws.Range("B14:F14").Value = "-"

# Row 15 (new, shifted from the old row 14): 16:40 slot
$ws.Range("A15").Value = "16:40"
$ws.Range("B15:F15").Value = "-"

# Row 16 (new): 17:30 slot
$ws.Range("A16").Value = "17:30"
$ws.Range("B16:F16").Value = "-"

# Row 17 (new): 18:20 slot, remaining columns left blank
$ws.Range("A17").Value = "18:20"
$ws.Range("B17:F17").Value = ""
